# Fruta / hortaliza, semanal
# Insert a new week of Kiwi price data (4 quality grades, "bandeja 10 kilos",
# Provincia de Curicó) at the top of the data table (rows 311-314), pushing
# the existing rows down by 4. The new rows mirror the constant columns
# (A-C, E-K) shared by every row in this single-commodity sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows right before the current row 311 (shifts 311:388 -> 315:392)
$ws.Rows("311:314").Insert()

# Values shared by every data row in this sheet (constant commodity columns)
$colA = 9
$colB = "Vega Central Mapocho de Santiago"
$colC = "Metropolitana"
$colE = 13
$colF = "Fruta"
$colG = 100101
$colH = "Berries"
$colI = 100101007
$colJ = "Kiwi"
$colK = "Hayward"

# Per-row specific values: Date, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm,
# Unidad, Origen, Precio $/Kg, Kg/unidad
$newRows = @(
    @{ Row=311; D=44476; L="Especial";                M=410; N=11000; O=11000; P=11000; Q="$/bandeja 10 kilos"; R="Provincia de Curicó"; S=1100; T=10 },
    @{ Row=312; D=44476; L="Extra (doble especial)";   M=400; N=12000; O=12000; P=12000; Q="$/bandeja 10 kilos"; R="Provincia de Curicó"; S=1200; T=10 },
    @{ Row=313; D=44476; L="Primera";                  M=350; N=9000;  O=9000;  P=9000;  Q="$/bandeja 10 kilos"; R="Provincia de Curicó"; S=900;  T=10 },
    @{ Row=314; D=44476; L="Segunda";                  M=300; N=7000;  O=7000;  P=7000;  Q="$/bandeja 10 kilos"; R="Provincia de Curicó"; S=700;  T=10 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $colA
    $ws.Cells.Item($row, 2).Value = $colB
    $ws.Cells.Item($row, 3).Value = $colC
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $colE
    $ws.Cells.Item($row, 6).Value = $colF
    $ws.Cells.Item($row, 7).Value = $colG
    $ws.Cells.Item($row, 8).Value = $colH
    $ws.Cells.Item($row, 9).Value = $colI
    $ws.Cells.Item($row, 10).Value = $colJ
    $ws.Cells.Item($row, 11).Value = $colK
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
